# Update the loan product "repaymentstrategy" scenario from "RBI (India)"
# to "Overdue/Due Fee/Int,Principal" on the ProductLoanInput sheet, and
# move the active selection to that cell (B17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

$ws.Activate()
$ws.Range("B17").Select()
